$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates - repulled data / recalculated final delta values
$updates = @{
    3  = -3
    5  = -5
    7  = -3
    10 = -9
    11 = -9
    12 = -6
    14 = -7
    15 = -2
    16 = -2
    17 = -3
    22 = 5
    23 = -5
    30 = 15
    32 = 10
    35 = -5
    36 = -4
    37 = -3
    38 = -9
    40 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
